# Generate Report for Handback
# Adds a second handback row (d1724f0c-...) alongside the existing
# (renamed) c9a46d89 -> 5cda55df row across the Overview / zh-cn / de-de
# sheets, expanding each sheet's table and refreshing the hyperlinks.

$wb = $excel.ActiveWorkbook

$oldGuid = "c9a46d89-45ba-4db9-838f-7659a4255cb9"
$newGuid1 = "5cda55df-95e1-436c-b019-5ed0b32e57b5"
$newGuid2 = "d1724f0c-1c19-4ace-ad47-e6b615a77776"

$zhHash1 = "8f26b7bd62c363764a0e1f4eb1f5c79889efba3e"
$zhHash2 = "99a40a0743b1ea5f069bea72f4fd569fdc670669"

# ---------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 2: rename the handback file from c9a46d89 -> 5cda55df
$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid1.md"

# Row 3: new handback file d1724f0c
$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("B3").Value = "e2e\$newGuid2.md"
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G3").Value = "2016-08-19 02:57:52"

# Rebuild hyperlinks for column B (display text must track the new names)
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/388d589be37742e13bca914744a031051fee7a61/e2e/$newGuid1.md", "", "", "e2e\$newGuid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/388d589be37742e13bca914744a031051fee7a61/e2e/$newGuid2.md", "", "", "e2e\$newGuid2.md") | Out-Null

# Expand the Overview table + autofilter to include row 3
$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G3")) | Out-Null

# ---------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 2: rename c9a46d89 -> 5cda55df (file + hash stay 8f26b7bd...)
$wsZh.Range("A2").Value = "$newGuid1.md"
$wsZh.Range("G2").Value = "$newGuid1.$zhHash1.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-19 02:57:46"
$wsZh.Range("I2").Value = "$newGuid1.md"
$wsZh.Range("J2").Value = "$newGuid1.$zhHash1.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-19 02:58:07"

# Row 3: new handback file d1724f0c
$wsZh.Range("A3").Value = "$newGuid2.md"
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$newGuid2.$zhHash2.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-19 02:57:46"
$wsZh.Range("I3").Value = "$newGuid2.md"
$wsZh.Range("J3").Value = "$newGuid2.$zhHash2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-19 02:58:07"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

# Rebuild hyperlinks for columns A and I
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/388d589be37742e13bca914744a031051fee7a61/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e1d3a7345a7e1511604fd6e81cace54be206908c/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/388d589be37742e13bca914744a031051fee7a61/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e1d3a7345a7e1511604fd6e81cace54be206908c/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null

# Expand the zh-cn table + autofilter to include row 3
$tblZh = $wsZh.ListObjects.Item(1)
$tblZh.Resize($wsZh.Range("A1:P3")) | Out-Null

# ---------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 2: rename c9a46d89 -> 5cda55df
$wsDe.Range("A2").Value = "$newGuid1.md"
$wsDe.Range("G2").Value = "$newGuid1.$zhHash1.de-de.xlf"
$wsDe.Range("I2").Value = "$newGuid1.md"
$wsDe.Range("J2").Value = "$newGuid1.$zhHash1.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-19 02:58:15"

# Row 3: new handback file d1724f0c
$wsDe.Range("A3").Value = "$newGuid2.md"
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$newGuid2.$zhHash2.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-19 02:57:52"
$wsDe.Range("I3").Value = "$newGuid2.md"
$wsDe.Range("J3").Value = "$newGuid2.$zhHash2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-19 02:58:15"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""

# Rebuild hyperlinks for columns A and I
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/388d589be37742e13bca914744a031051fee7a61/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6518937b4bfaa1ce206fa632967570f2f70f788f/e2e/$newGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/388d589be37742e13bca914744a031051fee7a61/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6518937b4bfaa1ce206fa632967570f2f70f788f/e2e/$newGuid2.md", "", "", "$newGuid2.md") | Out-Null

# Expand the de-de table + autofilter to include row 3
$tblDe = $wsDe.ListObjects.Item(1)
$tblDe.Resize($wsDe.Range("A1:P3")) | Out-Null

Write-Output "Handback report rows generated."
